{"js": "// Replace each \"old\u00f7divisor=\" arithmetic prompt in the worksheet table with\n// its new value, per the commit's regenerated problem set. Old -> new pairs\n// are all distinct strings, so a straightforward search/replace per pair is\n// safe and order-independent.\nconst replacements = [\n  [\"926\u00f77=\", \"640\u00f76=\"],\n  [\"486\u00f74=\", \"300\u00f72=\"],\n  [\"873\u00f79=\", \"419\u00f76=\"],\n  [\"696\u00f73=\", \"620\u00f78=\"],\n  [\"269\u00f78=\", \"496\u00f76=\"],\n  [\"179\u00f72=\", \"855\u00f77=\"],\n  [\"831\u00f77=\", \"195\u00f76=\"],\n  [\"723\u00f75=\", \"961\u00f73=\"],\n  [\"360\u00f72=\", \"915\u00f79=\"],\n  [\"987\u00f78=\", \"626\u00f75=\"],\n  [\"539\u00f77=\", \"472\u00f79=\"],\n  [\"778\u00f78=\", \"719\u00f76=\"],\n  [\"118\u00f78=\", \"894\u00f79=\"],\n  [\"619\u00f78=\", \"938\u00f73=\"],\n  [\"746\u00f76=\", \"757\u00f74=\"],\n  [\"397\u00f76=\", \"606\u00f76=\"],\n  [\"350\u00f78=\", \"579\u00f75=\"],\n  [\"222\u00f79=\", \"215\u00f77=\"],\n  [\"632\u00f73=\", \"129\u00f73=\"],\n  [\"884\u00f75=\", \"342\u00f72=\"],\n  [\"852\u00f72=\", \"107\u00f77=\"],\n  [\"650\u00f75=\", \"465\u00f76=\"],\n  [\"224\u00f77=\", \"634\u00f74=\"],\n  [\"980\u00f77=\", \"221\u00f77=\"],\n  [\"265\u00f72=\", \"555\u00f74=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each \"old\u00f7divisor=\" arithmetic prompt in the worksheet table with\n# its new value, per the commit's regenerated problem set.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"926\u00f77=\", \"640\u00f76=\"),\n  @(\"486\u00f74=\", \"300\u00f72=\"),\n  @(\"873\u00f79=\", \"419\u00f76=\"),\n  @(\"696\u00f73=\", \"620\u00f78=\"),\n  @(\"269\u00f78=\", \"496\u00f76=\"),\n  @(\"179\u00f72=\", \"855\u00f77=\"),\n  @(\"831\u00f77=\", \"195\u00f76=\"),\n  @(\"723\u00f75=\", \"961\u00f73=\"),\n  @(\"360\u00f72=\", \"915\u00f79=\"),\n  @(\"987\u00f78=\", \"626\u00f75=\"),\n  @(\"539\u00f77=\", \"472\u00f79=\"),\n  @(\"778\u00f78=\", \"719\u00f76=\"),\n  @(\"118\u00f78=\", \"894\u00f79=\"),\n  @(\"619\u00f78=\", \"938\u00f73=\"),\n  @(\"746\u00f76=\", \"757\u00f74=\"),\n  @(\"397\u00f76=\", \"606\u00f76=\"),\n  @(\"350\u00f78=\", \"579\u00f75=\"),\n  @(\"222\u00f79=\", \"215\u00f77=\"),\n  @(\"632\u00f73=\", \"129\u00f73=\"),\n  @(\"884\u00f75=\", \"342\u00f72=\"),\n  @(\"852\u00f72=\", \"107\u00f77=\"),\n  @(\"650\u00f75=\", \"465\u00f76=\"),\n  @(\"224\u00f77=\", \"634\u00f74=\"),\n  @(\"980\u00f77=\", \"221\u00f77=\"),\n  @(\"265\u00f72=\", \"555\u00f74=\"),\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $old\n  $find.Replacement.Text = $new\n  $find.Forward = $true\n  $find.Wrap = 1\n  $find.Execute($old, $false, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
